$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (font/wrap) from the existing B26 "response" cell onto the
# new B27 cell so it reuses the same cell style used throughout column B.
$ws.Range("B26").Copy() | Out-Null
$ws.Range("B27").PasteSpecial(-4122) | Out-Null

# New row content: tag/response pair for "thank-you"
$ws.Range("A27").Value = "thank-you"
$ws.Range("B27").Value = "ยินดีค่ะ ขอบคุณที่ใช้บริการกับเพจเรานะคะ โอกาสหน้าเชิญใหม่ค่ะ"

# Match the row height used by similarly-sized two-line rows (e.g. row 3).
$ws.Rows.Item(27).RowHeight = 25.5

# Update the sheet view so the new row is visible / selected, matching the
# author's saved cursor position after adding the row.
$ws.Range("C30").Select() | Out-Null
$excel.CutCopyMode = $false

Write-Host "Added thank-you row and response."
